{"js": "// Office.js (Word JavaScript API) script.\n// Applies the \"Reduce register wizard font sizes for cleaner balance\" edit:\n//  - Rename \"Register Page Wizard\" row to \"Register Wizard Typography\" and\n//    update its \"What's Done\" / \"What's Pending\" cells.\n//  - Remove the \"Header Menu\" row entirely.\n//  - Update \"Onboarding Flow Logic\" row's \"What's Done\" cell.\n//  - Bump the last-pushed-commit hash.\n//  - Update the final status line.\n\nconst body = context.document.body;\n\n// Simple text replacements located via body.search (exact, case-sensitive).\nasync function replaceOnce(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n\n// 1. Row title: \"Register Page Wizard\" -> \"Register Wizard Typography\"\nawait replaceOnce(\n  \"Register Page Wizard\",\n  \"Register Wizard Typography\"\n);\n\n// 2. \"What's Done\" cell for that row.\nawait replaceOnce(\n  \"`/auth/register` now uses the onboarding wizard layout and flow (same color/font style as latest wizard pass), matching your requested 5-step experience on the register URL.\",\n  \"Reduced oversized typography across `/auth/register` wizard (left hero, step rail, main title/subtitle, labels, and CTA button) for a more balanced professional scale.\"\n);\n\n// 3. \"What's Pending / Partial\" cell for that row.\nawait replaceOnce(\n  \"Can still fine-tune tiny spacing per your next screenshot if needed.\",\n  \"If needed, we can reduce one more step after reviewing latest live screenshot.\"\n);\n\n// 4. Delete the whole \"Header Menu\" paragraph.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet headerMenuParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Header Menu\\t\") === 0) {\n    headerMenuParagraph = paragraphs.items[i];\n    break;\n  }\n}\nif (!headerMenuParagraph) {\n  throw new Error(\"Could not locate the 'Header Menu' paragraph.\");\n}\nheaderMenuParagraph.delete();\nawait context.sync();\n\n// 5. \"Onboarding Flow Logic\" row's \"What's Done\" cell.\nawait replaceOnce(\n  \"Flow unchanged: Register -> Verify OTP(email+mobile) -> Choose Plan -> Payment(if paid) -> Setup Store -> Activate -> redirect admin.\",\n  \"No logic change; full 5-step flow remains intact.\"\n);\n\n// 6. Last pushed commit hash.\nawait replaceOnce(\n  \"- Last pushed commit: 75bfc21\",\n  \"- Last pushed commit: 11c9568\"\n);\n\n// 7. Final status line.\nawait replaceOnce(\n  \"- Current register wizard mapping + header-menu removal is local and pending push.\",\n  \"- Current font-size reduction pass is local and pending push.\"\n);\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# Applies the \"Reduce register wizard font sizes for cleaner balance\" edit:\n#  - Rename \"Register Page Wizard\" row to \"Register Wizard Typography\" and\n#    update its \"What's Done\" / \"What's Pending\" cells.\n#  - Remove the \"Header Menu\" row entirely.\n#  - Update \"Onboarding Flow Logic\" row's \"What's Done\" cell.\n#  - Bump the last-pushed-commit hash.\n#  - Update the final status line.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $rng = $d.Content\n    [void]$rng.Find.Execute(\n        $findText,      # FindText\n        $false,         # MatchCase\n        $false,         # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        1,              # Wrap (wdFindContinue)\n        $false,         # Format\n        $replaceText,   # ReplaceWith\n        2               # Replace (wdReplaceOne)\n    )\n}\n\n# 1. Row title: \"Register Page Wizard\" -> \"Register Wizard Typography\"\nReplace-Text \"Register Page Wizard\" \"Register Wizard Typography\"\n\n# 2. \"What's Done\" cell for that row.\nReplace-Text \"``/auth/register`` now uses the onboarding wizard layout and flow (same color/font style as latest wizard pass), matching your requested 5-step experience on the register URL.\" \"Reduced oversized typography across ``/auth/register`` wizard (left hero, step rail, main title/subtitle, labels, and CTA button) for a more balanced professional scale.\"\n\n# 3. \"What's Pending / Partial\" cell for that row.\nReplace-Text \"Can still fine-tune tiny spacing per your next screenshot if needed.\" \"If needed, we can reduce one more step after reviewing latest live screenshot.\"\n\n# 4. Delete the whole \"Header Menu\" paragraph (including its paragraph mark)\n#    without disturbing any other paragraph's run/tab structure.\n$targetParagraph = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.StartsWith(\"Header Menu`t\")) {\n        $targetParagraph = $p\n        break\n    }\n}\nif ($targetParagraph -eq $null) {\n    throw \"Could not locate the 'Header Menu' paragraph.\"\n}\n$targetParagraph.Range.Delete()\n\n# 5. \"Onboarding Flow Logic\" row's \"What's Done\" cell.\nReplace-Text \"Flow unchanged: Register -> Verify OTP(email+mobile) -> Choose Plan -> Payment(if paid) -> Setup Store -> Activate -> redirect admin.\" \"No logic change; full 5-step flow remains intact.\"\n\n# 6. Last pushed commit hash.\nReplace-Text \"- Last pushed commit: 75bfc21\" \"- Last pushed commit: 11c9568\"\n\n# 7. Final status line.\nReplace-Text \"- Current register wizard mapping + header-menu removal is local and pending push.\" \"- Current font-size reduction pass is local and pending push.\"\n"}
